$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell without letting Excel's
# type-inference turn numeric-looking strings (e.g. "11680391") into a
# number. We build the text in a scratch cell via a formula that evaluates
# to a string, copy it, and Paste-Special "values only" into the target -
# this preserves the string type and keeps the cell's existing style.
function Set-TextValue {
    param($targetAddr, [string]$text)

    $scratch = $ws.Range("Z100")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
    $scratch.Clear()
}

# --- Header row ---
Set-TextValue "B1" "Word Choice"

# --- Data rows: HashMethod (col A) stays the same; col B becomes the
# constant "Elucidate"; col C gets the new Web_Result hash value. ---
$webResults = @{
    2  = "bca0ac31500f40bdc1000e0aef880e73c817b5600fccf89dda1c198ed0b9e67e"
    3  = "11680391"
    4  = "ad93dfac"
    5  = "da809830"
    6  = "f8248dd2ec3b670a7a310407d0943c3f"
    7  = "7c8fc2ccdc68148861edf6549f782d3b2a4f7e0d"
    8  = "300bf3138202e8ac8d1c8fc67e40f476d29d01bf34a2f94e"
    9  = "77b84524ad0ab4fbf82c6e97aca950bb7d89919a5063521aa211267b"
    10 = "05740fa99fbd750d66dc32b4c9076d1aff28c927ed3f17533aff918b5d899c91"
    11 = "fc86045e7238f0466ad44662bb86f987"
    12 = "a4efc2a92c046b344f9a5e73cfd6f86492770b5c"
    13 = "f965b84a40c0f4af660918fcd088d5defe849cf7359b82eb"
    14 = "f6d26fe7f470da9c803bd3888b96e1eb165fb82bfda4207d640ce4ee"
    15 = "090bc3a2c1cf501784b6b2bdd81d872c52c146d27ecf45adef7246bec080c5bf"
    16 = "c7cf0c1bbc8b60f4767668fb0ebb7b80"
    17 = "2252ad7ccbf0a54ceb21d7a8f275a489fed6f773"
    18 = "ff7744673ce7c4ed52db4e83be6d7969215077905c8b77f2"
    19 = "889d3c5bb93249de49ff3ee1d7b1d0ddf46d833fbe3a73fd8d01271b"
    20 = "e498784f5dda5d5ff23d1814df87a7a77b868d27d2549cbac21d728f31e7be38"
}

for ($row = 2; $row -le 20; $row++) {
    Set-TextValue ("B" + $row) "Elucidate"
    Set-TextValue ("C" + $row) $webResults[$row]
}

# The old "Match" column (D) is gone in the new layout.
$ws.Columns.Item(4).Delete()

# Mirror the author's final selection (cell B2).
$ws.Range("B2").Select()
